# Update the "想去人数" (F column) counters that changed in the latest
# gh-pages data refresh (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 1204
$ws1.Range("F8").Value  = 288
$ws1.Range("F10").Value = 1244
$ws1.Range("F11").Value = 28078
$ws1.Range("F12").Value = 3258
$ws1.Range("F15").Value = 450
$ws1.Range("F22").Value = 242
$ws1.Range("F25").Value = 23
$ws1.Range("F28").Value = 86

# --- Sheet "演出" (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value  = 86
$ws2.Range("F11").Value = 4230
$ws2.Range("F22").Value = 4224

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F13").Value = 1204
$ws4.Range("F16").Value = 288
$ws4.Range("F20").Value = 86
$ws4.Range("F21").Value = 86
$ws4.Range("F30").Value = 450
$ws4.Range("F39").Value = 23
$ws4.Range("F43").Value = 86
